$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: coupon code
$ws.Range("A2").Value = "DV426SM5VGXG"

# B2: original value 100 -> 50
$ws.Range("B2").Value = 50

# C2: cost 0 -> 25
$ws.Range("C2").Value = 25

# D2: company
$ws.Range("D2").Value = "Wolt"

# E2: description
$ws.Range("E2").Value = "תו קנייה בשווי 50 ₪"

# F2: expiry date, keep as plain text (not an Excel date serial)
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2030-02-02"
$ws.Range("F2").Style = "Normal"

# G2 (תגיות) and H2 (סטטוס) are unchanged by the diff
